$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date value for rows 2-6 from 45174 to 45175
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
